$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.996.37"
$ws.Range("E2").Value = "  -1.98%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.902.94"
$ws.Range("E3").Value = "  -4.25%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.45"
$ws.Range("E5").Value = "  -0.98%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4596"
$ws.Range("E7").Value = "  -1.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3806"
$ws.Range("E8").Value = "  -2.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07703"
$ws.Range("E9").Value = "  -2.90%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9725"
$ws.Range("E10").Value = "  -2.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.88"
$ws.Range("E11").Value = "  -4.40%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.908.58"
$ws.Range("E12").Value = "  -3.89%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.919"
$ws.Range("E13").Value = "  -3.79%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.643"
$ws.Range("E14").Value = "  -3.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07086"
$ws.Range("E15").Value = "  -0.28%  "
$ws.Range("E16").Value = "  +0.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "83.65"
$ws.Range("E17").Value = "  -4.69%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000009473"
$ws.Range("E18").Value = "  -4.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.61"
$ws.Range("E19").Value = "  -4.16%  "
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "28.971.22"
$ws.Range("E21").Value = "  -2.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.295"
$ws.Range("E22").Value = "  -4.65%  "
$ws.Range("E23").Value = "  -2.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.099"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "157.88"
$ws.Range("E25").Value = "  -0.91%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "19.08"
$ws.Range("E26").Value = "  -2.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "5.599"
$ws.Range("E27").Value = "  -3.74%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "117.41"
$ws.Range("E28").Value = "  -2.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.846"
$ws.Range("E29").Value = "  -2.75%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09250"
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.8557"
$ws.Range("E31").Value = "  -4.40%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.084"
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.236"
$ws.Range("E33").Value = "  -7.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.953"
$ws.Range("E34").Value = "  -7.50%  "
$ws.Range("E35").Value = "  -2.32%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.138"
$ws.Range("E36").Value = "  -3.21%  "
$ws.Range("E37").Value = "  +0.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02031"
$ws.Range("E38").Value = "  -3.31%  "
$ws.Range("B39").Value = "FraxShare"
$ws.Range("C39").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.391"
$ws.Range("E39").Value = "  -5.77%  "
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5473"
$ws.Range("E40").Value = "  -4.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1751"
$ws.Range("E41").Value = "  -3.10%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "9.275"
$ws.Range("E42").Value = "  -4.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.758"
$ws.Range("E43").Value = "  -1.11%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5152"
$ws.Range("E44").Value = "  -3.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.19"
$ws.Range("E45").Value = "  -6.16%  "
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06840"
$ws.Range("E46").Value = "  -1.54%  "
$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.076"
$ws.Range("E47").Value = "  -4.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.000002581"
$ws.Range("E48").Value = "  -19.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "110.07"
$ws.Range("E49").Value = "  -3.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.765"
$ws.Range("E51").Value = "  +0.14%  "
